{"js": "// Replace the division problems in the table with their updated values.\n// Each old value is unique within the document, so a simple matchCase\n// search-and-replace (in document order) reproduces the diff exactly.\nconst replacements = [\n  [\"752\u00f74=\", \"770\u00f76=\"],\n  [\"147\u00f78=\", \"464\u00f76=\"],\n  [\"628\u00f75=\", \"714\u00f77=\"],\n  [\"776\u00f78=\", \"407\u00f77=\"],\n  [\"820\u00f79=\", \"931\u00f78=\"],\n  [\"554\u00f78=\", \"743\u00f72=\"],\n  [\"401\u00f74=\", \"621\u00f79=\"],\n  [\"619\u00f79=\", \"700\u00f77=\"],\n  [\"818\u00f73=\", \"155\u00f73=\"],\n  [\"258\u00f75=\", \"692\u00f78=\"],\n  [\"687\u00f72=\", \"864\u00f73=\"],\n  [\"694\u00f77=\", \"712\u00f72=\"],\n  [\"166\u00f74=\", \"291\u00f77=\"],\n  [\"605\u00f73=\", \"617\u00f72=\"],\n  [\"626\u00f75=\", \"699\u00f78=\"],\n  [\"974\u00f74=\", \"864\u00f73=\"],\n  [\"221\u00f73=\", \"578\u00f76=\"],\n  [\"182\u00f75=\", \"629\u00f76=\"],\n  [\"974\u00f77=\", \"105\u00f78=\"],\n  [\"702\u00f72=\", \"401\u00f77=\"],\n  [\"826\u00f77=\", \"306\u00f79=\"],\n  [\"532\u00f73=\", \"685\u00f73=\"],\n  [\"872\u00f75=\", \"582\u00f79=\"],\n  [\"587\u00f76=\", \"389\u00f79=\"],\n  [\"289\u00f79=\", \"856\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division problems in the table with their updated values.\n# Each \"old\" value is unique within the document, so Find/Replace (one\n# call per pair, in document order) reproduces the diff exactly.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"752\u00f74=\"; New = \"770\u00f76=\" },\n    @{ Old = \"147\u00f78=\"; New = \"464\u00f76=\" },\n    @{ Old = \"628\u00f75=\"; New = \"714\u00f77=\" },\n    @{ Old = \"776\u00f78=\"; New = \"407\u00f77=\" },\n    @{ Old = \"820\u00f79=\"; New = \"931\u00f78=\" },\n    @{ Old = \"554\u00f78=\"; New = \"743\u00f72=\" },\n    @{ Old = \"401\u00f74=\"; New = \"621\u00f79=\" },\n    @{ Old = \"619\u00f79=\"; New = \"700\u00f77=\" },\n    @{ Old = \"818\u00f73=\"; New = \"155\u00f73=\" },\n    @{ Old = \"258\u00f75=\"; New = \"692\u00f78=\" },\n    @{ Old = \"687\u00f72=\"; New = \"864\u00f73=\" },\n    @{ Old = \"694\u00f77=\"; New = \"712\u00f72=\" },\n    @{ Old = \"166\u00f74=\"; New = \"291\u00f77=\" },\n    @{ Old = \"605\u00f73=\"; New = \"617\u00f72=\" },\n    @{ Old = \"626\u00f75=\"; New = \"699\u00f78=\" },\n    @{ Old = \"974\u00f74=\"; New = \"864\u00f73=\" },\n    @{ Old = \"221\u00f73=\"; New = \"578\u00f76=\" },\n    @{ Old = \"182\u00f75=\"; New = \"629\u00f76=\" },\n    @{ Old = \"974\u00f77=\"; New = \"105\u00f78=\" },\n    @{ Old = \"702\u00f72=\"; New = \"401\u00f77=\" },\n    @{ Old = \"826\u00f77=\"; New = \"306\u00f79=\" },\n    @{ Old = \"532\u00f73=\"; New = \"685\u00f73=\" },\n    @{ Old = \"872\u00f75=\"; New = \"582\u00f79=\" },\n    @{ Old = \"587\u00f76=\"; New = \"389\u00f79=\" },\n    @{ Old = \"289\u00f79=\"; New = \"856\u00f76=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $r.Old,   # FindText\n        $true,    # MatchCase\n        $false,   # MatchWholeWord\n        $false,   # MatchWildcards\n        $false,   # MatchSoundsLike\n        $false,   # MatchAllWordForms\n        $true,    # Forward\n        0,        # Wrap (wdFindStop)\n        $false,   # Format\n        $r.New,   # ReplaceWith\n        2         # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n\n$d.Save()\n"}
